$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 (Q0-Q3) with new values
$ws.Range("B2").Value = 0.9727825379741304
$ws.Range("C2").Value = 1.078664240660259
$ws.Range("D2").Value = 1.922255377816953
$ws.Range("E2").Value = 1.386454246564578
$ws.Range("F2").Value = 1.02519390896829
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.8757695348867685
$ws.Range("C3").Value = 0.9186024224941243
$ws.Range("D3").Value = 1.56012335503454
$ws.Range("E3").Value = 1.249048980238381
$ws.Range("F3").Value = 0.9387646000040889
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = 1.137903216480779
$ws.Range("C4").Value = 1.137903216480779
$ws.Range("D4").Value = 2.3596763441606
$ws.Range("E4").Value = 1.536123804958637
$ws.Range("F4").Value = 1.130408393855935
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = 1.25991956803407
$ws.Range("C5").Value = 1.25991956803407
$ws.Range("D5").Value = 2.455580339113391
$ws.Range("E5").Value = 1.56702914430887
$ws.Range("F5").Value = 1.317712427806791
$ws.Range("G5").Value = 2

# Remove rows 6-9 (Q4, Q5, Q6, Q7) entirely
$ws.Range("A6:G9").EntireRow.Delete()
